$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the indicator name text in B4 (17.4.1 indicator description changed)
$ws.Range("B4").Value = "17.4.1 Доля поступлений от экспорта товаров и услуг и первичного дохода, расходуемая на обслуживание долга"

# Move/update the active selection to B4 to match the saved view state
$ws.Range("B4").Select()
